# Rename the "NewReferences" sheet to "New references".
# Renaming automatically updates the sheet's defined-name reference
# (_xlnm._FilterDatabase) to use the quoted sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewReferences")
$ws.Name = "New references"

# Make the renamed sheet the active tab (was previously "Tests").
$ws.Activate()
